# Resolved issue with exit button always taking to motorcycle menu
#
# Adds six new rows to the "Table1" listing on Sheet1 describing the
# in-call popup and active-call button behaviour (answer / reject /
# volume / exit), plus a trailing blank formatted row, and extends the
# table/autofilter/dimension to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 150 is the current last data row of the table and carries the
# "bottom of table" border formatting (s=7/8/9). Stamp that same
# formatting onto every new row (151-157) first, then fill in values.
for ($i = 151; $i -le 157; $i++) {
    $ws.Range("A150:D150").Copy()
    $ws.Range("A" + $i + ":D" + $i).PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = $false

# Fill in the new data. Row 154 (DOWN) / 156 (EXIT) are populated before
# their siblings 153 (UP) / 155 (SET) so that newly introduced strings
# land in the workbook's shared string table in the same order as the
# original authoring session.
$ws.Range("A151").Value = 150
$ws.Range("B151").Value = "POPUP_INCALL_INSTRUCTION"
$ws.Range("C151").Value = "SET"
$ws.Range("D151").Value = "BLU_FnIncomingCallAnswer()"

$ws.Range("A152").Value = 151
$ws.Range("B152").Value = "POPUP_INCALL_INSTRUCTION"
$ws.Range("C152").Value = "EXIT"
$ws.Range("D152").Value = "BLU_FnIncomingCallReject()"

$ws.Range("A154").Value = 153
$ws.Range("B154").Value = "ACTIVE_CALL"
$ws.Range("C154").Value = "DOWN"
$ws.Range("D154").Value = "BLU_FnCallVolDec()"

$ws.Range("A153").Value = 152
$ws.Range("B153").Value = "ACTIVE_CALL"
$ws.Range("C153").Value = "UP"
$ws.Range("D153").Value = "BLU_FnCallVolInc()"

$ws.Range("A156").Value = 155
$ws.Range("B156").Value = "ACTIVE_CALL"
$ws.Range("C156").Value = "EXIT"
$ws.Range("D156").Value = "BLU_FnActvCallEXIT()"

$ws.Range("A155").Value = 154
$ws.Range("B155").Value = "ACTIVE_CALL"
$ws.Range("C155").Value = "SET"
$ws.Range("D155").Value = "BLU_FnActvCallSET()"

# Row 157 stays blank (only formatting carried over above).

# Resize the table / autofilter range so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D157"))

# Scroll down and select the newly added Screen/Action/Function cells,
# matching where the author was working when the file was saved.
$ws.Application.Goto($ws.Range("A139"), $false)
$ws.Range("B151:D156").Select()
